# "rozmery po uprave velikosti boxu" - add new device rows with updated
# dimensions after resizing the box, on List1 (sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11: Nexus 6P, landscape - resolution only (entered first so the
# shared-string pool interns "Nexus 6P" before "DP")
$ws.Range("A11").Value = "Nexus 6P"
$ws.Range("B11").Value = "landscape"
$ws.Range("C11").Value = 2560
$ws.Range("D11").Value = 1440

# Row 10: DP - only E/F given (box size), J10 = E10/F10
$ws.Range("A10").Value = "DP"
$ws.Range("E10").Value = 120
$ws.Range("F10").Value = 180
$ws.Range("J10").Formula = "=E10/F10"

# Row 12: LeliMath, landscape - full row incl. box + letter size
$ws.Range("A12").Value = "LeliMath"
$ws.Range("B12").Value = "landscape"
$ws.Range("C12").Value = 2390
$ws.Range("D12").Value = 1190
$ws.Range("E12").Value = 420
$ws.Range("F12").Value = 630
$ws.Range("G12").Value = 29
$ws.Range("H12").Value = 37
$ws.Range("J12").Formula = "=E12/F12"

# Row 13: Nexus 5x - resolution only
$ws.Range("A13").Value = "Nexus 5x"
$ws.Range("C13").Value = 1920
$ws.Range("D13").Value = 1080

# Row 14: LeliMath, portrait - full row incl. box + letter size
$ws.Range("A14").Value = "LeliMath"
$ws.Range("B14").Value = "portrait"
$ws.Range("C14").Value = 1080
$ws.Range("D14").Value = 1590
$ws.Range("E14").Value = 310
$ws.Range("F14").Value = 470
$ws.Range("G14").Value = 24
$ws.Range("H14").Value = 30
$ws.Range("J14").Formula = "=E14/F14"

# Page setup - A4, portrait
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave selection where the user finished typing (row below last entry)
$ws.Range("J15").Select() | Out-Null
